# Add a new payment row (row 44) to the "Feuil1" sheet, duplicating the
# donor "Nathalie Parasol" payment already recorded on row 41 but with a
# new transfer reference ("ESP5") and a different amount (15 instead of 20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 already carries the exact cell formatting (date format, wrapped
# text, postal-code-as-text, ...) that the new row needs, so clone its
# formatting first and then overwrite the cell values. This avoids
# creating any brand-new cell styles (the style table is unchanged by the
# real edit).
$ws.Range("A41:K41").Copy()
$ws.Range("A44:K44").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A44").Value = 45314
$ws.Range("B44").Value = "ESP5"
$ws.Range("C44").Value = "apn514@mail.com"
$ws.Range("D44").Value = "Parasol"
$ws.Range("E44").Value = "Nathalie"
$ws.Range("F44").Value = "39 rue de la gerbille"
$ws.Range("G44").Value = "38000"
$ws.Range("H44").Value = "GRENOBLE"
$ws.Range("I44").Value = 15
$ws.Range("J44").Value = "E"
$ws.Range("K44").Value = "N"

# Match the saved selection/cursor position recorded in the workbook.
$ws.Range("I45").Select() | Out-Null
